$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the Heading1 title
#    (paragraph 1). It has an empty leading run, a bold "Meta description"
#    run, and a plain run with the rest of the description text.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$metaXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Get ready to play Amazon Gold, a jungle-themed slot game with potential for high rewards and exciting features. Play for free now!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$metaPara.Range.InsertXML($metaXml) | Out-Null

# ---------------------------------------------------------------------------
# 2) Remove the trailing duplicate bold title paragraph near the end of the
#    document, and replace the italic "Get ready..." paragraph's text with
#    the new image-prompt text (keeping its italic formatting + leading
#    empty run untouched).
# ---------------------------------------------------------------------------
$boldTarget = "Play Amazon Gold for Free - Exciting Jungle-Themed Slot Game"
$boldTitlePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.Trim() -eq $boldTarget -and $i -ne 1) {
        $boldTitlePara = $cand
    }
}
$boldTitlePara.Range.Delete()

$italicTarget = "Get ready to play Amazon Gold, a jungle-themed slot game with potential for high rewards and exciting features. Play for free now!"
$imgPromptPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.Trim() -eq $italicTarget) {
        $imgPromptPara = $cand
    }
}

$pr = $imgPromptPara.Range
$textRange = $d.Range($pr.Start, $pr.End - 1)

$imgXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:i/></w:rPr><w:t>Create a cartoon-style feature image for Amazon Gold that showcases a happy Maya warrior wearing glasses. The warrior should be surrounded by the colorful and exotic flora and fauna of the jungle, with symbols from the game, such as parrots, crocodiles, monkeys, and gold nuggets, prominently featured. The image should also convey the excitement and potential for high rewards with the use of Free Spins and Wild multipliers. Make sure to include the game's title "Amazon Gold" in the image.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$textRange.InsertXML($imgXml) | Out-Null

Write-Host "Edits applied. Paragraph count:" $d.Paragraphs.Count
